$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
# Row 8 (Leve Item ID 4565)
$ws.Range("H8").Value = 108
$ws.Range("I8").Value = 108
$ws.Range("K8").Value = 324
$ws.Range("M8").Value = -185

# Row 48 (Leve Item ID 4587)
$ws.Range("H48").Value = 4756.75
$ws.Range("I48").Value = 3008.5
$ws.Range("J48").Value = 6505
$ws.Range("K48").Value = 9025.5
$ws.Range("L48").Value = 19515
$ws.Range("M48").Value = -8733.5
$ws.Range("N48").Value = -20099

# Row 49 (Leve Item ID 4588)
$ws.Range("H49").Value = 1230.8889
$ws.Range("I49").Value = 959.6667
$ws.Range("J49").Value = 1502.1111
$ws.Range("K49").Value = 2879.0001
$ws.Range("L49").Value = 4506.3333
$ws.Range("M49").Value = -2743.0001
$ws.Range("N49").Value = -4778.3333

# Row 56 (Leve Item ID 4587)
$ws.Range("H56").Value = 4756.75
$ws.Range("I56").Value = 3008.5
$ws.Range("J56").Value = 6505
$ws.Range("K56").Value = 9025.5
$ws.Range("L56").Value = 19515
$ws.Range("M56").Value = -8491.5
$ws.Range("N56").Value = -20583

# Row 58 (Leve Item ID 4606)
$ws.Range("H58").Value = 1465.3125
$ws.Range("I58").Value = 189
$ws.Range("J58").Value = 2045.4546
$ws.Range("K58").Value = 567
$ws.Range("L58").Value = 6136.3638
$ws.Range("M58").Value = -417
$ws.Range("N58").Value = -6436.3638

# Row 59 (Leve Item ID 4586)
$ws.Range("H59").Value = 977
$ws.Range("J59").Value = 977
$ws.Range("L59").Value = 2931
$ws.Range("N59").Value = -4045

# Row 61 (Leve Item ID 4604)
$ws.Range("H61").Value = 363.2
$ws.Range("I61").Value = 224.75
$ws.Range("J61").Value = 917
$ws.Range("K61").Value = 674.25
$ws.Range("L61").Value = 2751
$ws.Range("M61").Value = -502.25
$ws.Range("N61").Value = -3095

# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 730
$ws.Range("I98").Value = 715
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 715
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 783
$ws.Range("N98").Value = -3996

# Row 113 (Leve Item ID 27775)
$ws.Range("H113").Value = 2792.762
$ws.Range("I113").Value = 2400.8333
$ws.Range("J113").Value = 2949.5334
$ws.Range("K113").Value = 2400.8333
$ws.Range("L113").Value = 2949.5334
$ws.Range("M113").Value = 853.1667000000002
$ws.Range("N113").Value = -9457.5334

# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 730
$ws.Range("I122").Value = 715
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2145
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = 305
$ws.Range("N122").Value = -7900

# Row 125 (Leve Item ID 36228)
$ws.Range("H125").Value = 350.58334
$ws.Range("I125").Value = 266.7143
$ws.Range("J125").Value = 468
$ws.Range("K125").Value = 2400.4287
$ws.Range("L125").Value = 4212
$ws.Range("M125").Value = 59.57130000000006
$ws.Range("N125").Value = -9132

# Row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 839.1
$ws.Range("I129").Value = 292.6
$ws.Range("J129").Value = 1112.35
$ws.Range("K129").Value = 877.8000000000001
$ws.Range("L129").Value = 3337.05
$ws.Range("M129").Value = 4122.2
$ws.Range("N129").Value = -13337.05

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 1602.6508
$ws.Range("I138").Value = 568.11365
$ws.Range("J138").Value = 3998.4211
$ws.Range("K138").Value = 1704.34095
$ws.Range("L138").Value = 11995.2633
$ws.Range("M138").Value = 3435.65905
$ws.Range("N138").Value = -22275.2633

$ws = $wb.Sheets.Item("CRP")
# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 885.3594000000001
$ws.Range("I58").Value = 457.5
$ws.Range("J58").Value = 1370.2667
$ws.Range("K58").Value = 457.5
$ws.Range("L58").Value = 1370.2667
$ws.Range("M58").Value = -254.5
$ws.Range("N58").Value = -1776.2667

# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 5110209
$ws.Range("I99").Value = 5960577
$ws.Range("J99").Value = 8000
$ws.Range("K99").Value = 5960577
$ws.Range("L99").Value = 8000
$ws.Range("M99").Value = -5959079
$ws.Range("N99").Value = -10996

# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 5110209
$ws.Range("I126").Value = 5960577
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 17881731
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -17879261
$ws.Range("N126").Value = -28940

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 1832.8125
$ws.Range("I132").Value = 1690.1111
$ws.Range("J132").Value = 2016.2858
$ws.Range("K132").Value = 5070.3333
$ws.Range("L132").Value = 6048.857400000001
$ws.Range("M132").Value = -2540.3333
$ws.Range("N132").Value = -11108.8574

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 885.3594000000001
$ws.Range("I136").Value = 457.5
$ws.Range("J136").Value = 1370.2667
$ws.Range("K136").Value = 1372.5
$ws.Range("L136").Value = 4110.800099999999
$ws.Range("M136").Value = 1177.5
$ws.Range("N136").Value = -9210.8001

$ws = $wb.Sheets.Item("CUL")
# Row 80 (Leve Item ID 12890)
$ws.Range("H80").Value = 2000
$ws.Range("J80").Value = 2000
$ws.Range("L80").Value = 6000
$ws.Range("N80").Value = -7872

# Row 83 (Leve Item ID 12890)
$ws.Range("H83").Value = 2000
$ws.Range("J83").Value = 2000
$ws.Range("L83").Value = 18000
$ws.Range("N83").Value = -27360

# Row 114 (Leve Item ID 27865)
$ws.Range("H114").Value = 895.63336
$ws.Range("I114").Value = 1384.7693
$ws.Range("J114").Value = 521.58826
$ws.Range("K114").Value = 4154.3079
$ws.Range("L114").Value = 1564.76478
$ws.Range("M114").Value = -900.3078999999998
$ws.Range("N114").Value = -8072.76478

# Row 117 (Leve Item ID 27870)
$ws.Range("H117").Value = 1415.409
$ws.Range("I117").Value = 684.5
$ws.Range("J117").Value = 1689.5
$ws.Range("K117").Value = 2053.5
$ws.Range("L117").Value = 5068.5
$ws.Range("M117").Value = 1388.5
$ws.Range("N117").Value = -11952.5

# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 1021.4474
$ws.Range("J131").Value = 1059.2572
$ws.Range("L131").Value = 3177.7716
$ws.Range("N131").Value = -13257.7716

$ws = $wb.Sheets.Item("GSM")
# Row 107 (Leve Item ID 27802)
$ws.Range("H107").Value = 912.6111
$ws.Range("I107").Value = 961.2857
$ws.Range("J107").Value = 881.63635
$ws.Range("K107").Value = 961.2857
$ws.Range("L107").Value = 881.63635
$ws.Range("M107").Value = 958.7143
$ws.Range("N107").Value = -4721.63635

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 2614.48
$ws.Range("I126").Value = 1790.5
$ws.Range("J126").Value = 4079.3333
$ws.Range("K126").Value = 5371.5
$ws.Range("L126").Value = 12237.9999
$ws.Range("M126").Value = -2901.5
$ws.Range("N126").Value = -17177.9999

# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 1718.12
$ws.Range("I132").Value = 1638.6666
$ws.Range("J132").Value = 1791.4615
$ws.Range("K132").Value = 4915.9998
$ws.Range("L132").Value = 5374.3845
$ws.Range("M132").Value = -2385.9998
$ws.Range("N132").Value = -10434.3845

$ws = $wb.Sheets.Item("LTW")
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 1243.6666
$ws.Range("I40").Value = 1141.5714
$ws.Range("J40").Value = 1601
$ws.Range("K40").Value = 1141.5714
$ws.Range("L40").Value = 1601
$ws.Range("M40").Value = -1005.5714
$ws.Range("N40").Value = -1873

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 11579618
$ws.Range("I132").Value = 36777988
$ws.Range("J132").Value = 1989.6757
$ws.Range("K132").Value = 110333964
$ws.Range("L132").Value = 5969.0271
$ws.Range("M132").Value = -110331434
$ws.Range("N132").Value = -11029.0271

$ws = $wb.Sheets.Item("WVR")
# Row 33 (Leve Item ID 2734)
$ws.Range("H33").Value = 27000
$ws.Range("J33").Value = 27000
$ws.Range("L33").Value = 27000
$ws.Range("N33").Value = -27500

# Row 36 (Leve Item ID 2734)
$ws.Range("H36").Value = 27000
$ws.Range("J36").Value = 27000
$ws.Range("L36").Value = 27000
$ws.Range("N36").Value = -27500

# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 58826170
$ws.Range("I81").Value = 90911640
$ws.Range("J81").Value = 2801
$ws.Range("K81").Value = 181823280
$ws.Range("L81").Value = 5602
$ws.Range("M81").Value = -181822219
$ws.Range("N81").Value = -7724

# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 58826170
$ws.Range("I84").Value = 90911640
$ws.Range("J84").Value = 2801
$ws.Range("K84").Value = 909116400
$ws.Range("L84").Value = 28010
$ws.Range("M84").Value = -909111096
$ws.Range("N84").Value = -38618

# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 2029.4595
$ws.Range("I126").Value = 2065.3125
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 6195.9375
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -3725.9375
$ws.Range("N126").Value = -10340

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1290.4651
$ws.Range("I132").Value = 1022.3043
$ws.Range("K132").Value = 3066.9129
$ws.Range("M132").Value = -536.9129000000003
